$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the single-line mailing address "929 Story Road, San Jose CA 95122"
#    into two paragraphs: "929 Story Road" and "San Jose, CA 95122", keeping the
#    same paragraph/run formatting (Arial, sz 22) on both.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "929 Story Road, San Jose CA 95122`r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Clone the paragraph (including its mark/formatting) right after itself so the
    # new paragraph inherits identical pPr/rPr, then fill in the split text.
    $dup = $target.Range.Duplicate
    $dup.InsertParagraphAfter()

    $paragraphs = @($d.Paragraphs)
    $idx = -1
    for ($i = 0; $i -lt $paragraphs.Count; $i++) {
        if ($paragraphs[$i].Range.Text -eq "929 Story Road, San Jose CA 95122`r") {
            $idx = $i
            break
        }
    }

    $p1 = $paragraphs[$idx]
    $p2 = $paragraphs[$idx + 1]

    $r1 = $p1.Range.Duplicate
    $r1.MoveEnd(1, -1) | Out-Null
    $r1.Text = "929 Story Road"

    $r2 = $p2.Range.Duplicate
    $r2.MoveEnd(1, -1) | Out-Null
    $r2.Text = "San Jose, CA 95122"
}

# 3. Remove the empty "NoSpacing" paragraph that immediately follows the
#    "Board of Directors" line.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Vietnam Town Condominium Owners Association Board of Directors`r") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}
